# Remove every <w:contextualSpacing .../> element from the paragraph
# properties (w:pPr) throughout the document body. The commit simply
# dropped this element from each paragraph's formatting (regardless of
# its w:val), leaving everything else untouched.

$d = $word.ActiveDocument

# Pull the full package OOXML for the whole document's content range so
# the edit is performed as a single atomic pass (this keeps relationship
# ids - e.g. hyperlink rIds - intact, unlike doing this per paragraph).
$full = $d.Content
$pkg = $full.WordOpenXML

# Strip every <w:contextualSpacing/> (self-closing, with or without
# attributes such as w:val="0"/"1") element wherever it occurs.
$pkg = $pkg -replace '<w:contextualSpacing(\s[^>]*)?/>', ''

# Write the modified package XML back into the document (discard the
# return value so it isn't echoed to the output stream).
$null = $full.InsertXML($pkg)
